$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "Averaged Data" placeholder values ("N/A") across the
# TPS / Energy / Nakamoto / % of nodes / Strengths / Weaknesses
# columns for every consensus-mechanism row.
$ws.Range("B2:G10").Value = "N/A"

# Apply number formats appropriate to each column's data type, same
# as selecting the table column's data and changing its format.
$ws.Range("B2:B10").NumberFormat = "0.00"   # TPS
$ws.Range("C2:C10").NumberFormat = "0.00"   # Energy Use per Transaction
$ws.Range("D2:D10").NumberFormat = "0.00"   # Nakamoto Coefficient
$ws.Range("E2:E10").Style = "Percent"       # % of nodes required to take over network
$ws.Range("F2:F10").NumberFormat = "@"      # Strengths
$ws.Range("G2:G10").NumberFormat = "@"      # Weaknesses

# Restore the last active selection as recorded in the saved session.
$ws.Range("J16").Select()
